# Generate Report for Handoff
#
# For the six rows that correspond to files that just got handed off
# (2358dca1, 3d23d051, 4fd11d5e, 593673e7, 5bd08ba7, eba68848 -> rows
# 7,8,9,10,11,14 on every sheet), refresh the handoff bookkeeping:
#   - Overview!G  (Latest HO Xliff Generate Date)
#   - de-de!H     (Latest Handoff Datetime)
#   - zh-cn!H     (Latest Handoff Datetime)
#   - zh-cn!E and de-de!E (Priority) get marked "ht" (handoff type)

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 14)

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-25 14:19:49"

    $dede.Cells.Item($r, 8).Value = "2016-08-25 14:19:49"

    $zhcn.Cells.Item($r, 8).Value = "2016-08-25 14:19:44"

    $zhcn.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 5).Value = "ht"
}
